# Update "want to go" counts (column F) on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5446
$wsExhibit.Range("F9").Value = 518

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5446
$wsAll.Range("F11").Value = 518
